$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.162.29"
$ws.Range("E2").Value = "'  +1.55%  "
$ws.Range("D3").Value = "'1.794.02"
$ws.Range("E3").Value = "'  +1.82%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  -0.18%  "
$ws.Range("D5").Value = "'323.83"
$ws.Range("E5").Value = "'  -0.98%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  +0.13%  "
$ws.Range("D7").Value = "'0.4290"
$ws.Range("E7").Value = "'  -2.79%  "
$ws.Range("D8").Value = "'0.3630"
$ws.Range("E8").Value = "'  -3.16%  "
$ws.Range("D9").Value = "'44.73"
$ws.Range("E9").Value = "'  -1.45%  "
$ws.Range("D10").Value = "'0.07535"
$ws.Range("E10").Value = "'  -3.44%  "
$ws.Range("D11").Value = "'1.117"
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "'  -0.02%  "
$ws.Range("E13").Value = "'  -0.54%  "
$ws.Range("D14").Value = "'6.156"
$ws.Range("E14").Value = "'  -0.69%  "
$ws.Range("D15").Value = "'7.333"
$ws.Range("E15").Value = "'  -0.89%  "
$ws.Range("D16").Value = "'1.810.17"
$ws.Range("E16").Value = "'  +2.95%  "
$ws.Range("D17").Value = "'92.08"
$ws.Range("E17").Value = "'  +1.18%  "
$ws.Range("D18").Value = "'0.00001072"
$ws.Range("E18").Value = "'  -1.11%  "
$ws.Range("D19").Value = "'0.06369"
$ws.Range("E20").Value = "'  +0.05%  "
$ws.Range("D21").Value = "'17.20"
$ws.Range("E21").Value = "'  -0.96%  "
$ws.Range("D22").Value = "'5.984"
$ws.Range("E22").Value = "'  -3.39%  "
$ws.Range("D23").Value = "'28.155.58"
$ws.Range("E23").Value = "'  +1.37%  "
$ws.Range("E24").Value = "'  -2.48%  "
$ws.Range("D25").Value = "'2.171"
$ws.Range("E25").Value = "'  -6.18%  "
$ws.Range("D26").Value = "'159.35"
$ws.Range("E26").Value = "'  +3.55%  "
$ws.Range("D27").Value = "'20.38"
$ws.Range("E27").Value = "'  -2.40%  "
$ws.Range("D28").Value = "'2.012.89"
$ws.Range("E28").Value = "'  +2.87%  "
$ws.Range("D29").Value = "'2.232"
$ws.Range("E29").Value = "'  -5.59%  "
$ws.Range("D30").Value = "'127.96"
$ws.Range("E30").Value = "'  -0.95%  "
$ws.Range("D31").Value = "'1.170"
$ws.Range("E31").Value = "'  -3.96%  "
$ws.Range("D32").Value = "'5.819"
$ws.Range("E32").Value = "'  +0.40%  "
$ws.Range("D33").Value = "'0.08992"
$ws.Range("E33").Value = "'  -3.28%  "
$ws.Range("D34").Value = "'3.533"
$ws.Range("E34").Value = "'  -3.05%  "
$ws.Range("D35").Value = "'12.74"
$ws.Range("E35").Value = "'  +0.10%  "
$ws.Range("D36").Value = "'0.02360"

# Row 37/38: TheSandbox and InternetComputer(DFINITY) swap places, with updated Price/Volume values
$ws.Range("B37").Value = "'InternetComputer(DFINITY)"
$ws.Range("C37").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.095"
$ws.Range("E37").Value = "'  -0.20%  "
$ws.Range("B38").Value = "'TheSandbox"
$ws.Range("C38").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "'0.6486"
$ws.Range("E38").Value = "'  -0.46%  "

$ws.Range("D39").Value = "'0.06112"
$ws.Range("E39").Value = "'  -0.84%  "
$ws.Range("D40").Value = "'0.2118"
$ws.Range("E40").Value = "'  -3.22%  "
$ws.Range("E41").Value = "'  -0.88%  "
$ws.Range("D42").Value = "'1.428"
$ws.Range("E42").Value = "'  +0.47%  "
$ws.Range("D43").Value = "'7.949"
$ws.Range("E43").Value = "'  -0.93%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "'  +0.07%  "
$ws.Range("D45").Value = "'13.55"
$ws.Range("E45").Value = "'  -2.59%  "
$ws.Range("D46").Value = "'0.6016"
$ws.Range("E46").Value = "'  -0.27%  "
$ws.Range("E47").Value = "'  -1.60%  "
$ws.Range("D48").Value = "'125.28"
$ws.Range("E48").Value = "'  -0.81%  "
$ws.Range("D49").Value = "'1.994"
$ws.Range("E49").Value = "'  -0.44%  "
$ws.Range("D50").Value = "'1.152"
$ws.Range("E50").Value = "'  +0.49%  "
$ws.Range("D51").Value = "'0.06967"
$ws.Range("E51").Value = "'  +0.79%  "
